# "cool changes for baris form editable now"
# - rows 4..6 removed from Sheet1 (only the first two data rows remain)
# - remaining data rows (2 & 3) have their numeric metrics reset to 0
# - the free-text "Summoner Name" cell (G) is cleared out
# - the "Summoner Role" cell (H) is normalized to "SOLO" for both rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused data rows (4, 5, 6) entirely so the sheet shrinks to A1:K3.
$ws.Range("A4:K6").EntireRow.Delete()

# Row 2 -> everything zeroed out / blanked, role normalized to SOLO.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "SOLO"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Row 3 -> same treatment.
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "SOLO"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
